$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update last_edited_time (column D) for rows 2, 3, 6, 8, 11, 13
$rowsToStamp = @(2, 3, 6, 8, 11, 13)
foreach ($r in $rowsToStamp) {
    $ws.Range("D$r").Value = "2024-07-08T01:58:00.000Z"
}

# Update numeric property values on row 13 (report co so)
$ws.Range("S13").Value = 37717000
$ws.Range("W13").Value = 48283000
$ws.Range("AA13").Value = 3800000
$ws.Range("AE13").Value = 86000000
$ws.Range("AH13").Value = 70000000
$ws.Range("AK13").Value = 8
$ws.Range("AN13").Value = 16000000
$ws.Range("AQ13").Value = 73800000
